# Re-run of the face-recognition scoring pipeline: a new sample (FY_4.png)
# was inserted after FY_3.png, shifting the former TO_*/TD_* rows down by
# one, and the distance/threshold/score columns (B:E) were recomputed for
# every row. This writes the full refreshed A1:G45 table in one pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('AK_1.png', 0.361, 0.001, 0.496, 0.857, 'Akhlak Kamiswara', 'Benar'),
    @('AK_2.png', 0.694, 0.002, 0.487, 1, 'Akhlak Kamiswara', 'Benar'),
    @('AK_3.png', 0.668, 0.002, 0.649, 1, 'Akhlak Kamiswara', 'Benar'),
    @('AK_4.png', 0.335, 0.001, 0.47, 1, 'Akhlak Kamiswara', 'Benar'),
    @('AK_5.png', 0.319, 0.001, 0.551, 1, 'Akhlak Kamiswara', 'Benar'),
    @('MIB_1.png', 1.4, 0.005, 0.133, 0.714, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_2.png', 0.988, 0.003, 0.344, 1, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_3.png', 1.619, 0.005, 0.701, 1, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_4.png', 1.002, 0.003, 0.251, 1, 'Muhammad Iqbal Baqi', 'Benar'),
    @('MIB_5.png', 1.395, 0.005, 0.266, 1, 'Muhammad Iqbal Baqi', 'Benar'),
    @('AAH_1.png', 0.679, 0.002, 0.287, 1, 'Andrea Ayunove Hutami', 'Benar'),
    @('AAH_2.png', 0.944, 0.003, 0.721, 1, 'Andrea Ayunove Hutami', 'Benar'),
    @('AAH_3.png', 0.72, 0.002, 0.311, 1, 'Andrea Ayunove Hutami', 'Benar'),
    @('TI_1.png', 0.8149999999999999, 0.003, 0.502, 1, 'Toni Ismail', 'Benar'),
    @('TI_2.png', 0.786, 0.003, 0.498, 1, 'Toni Ismail', 'Benar'),
    @('TI_3.png', 0.5659999999999999, 0.002, 0.795, 1, 'Toni Ismail', 'Benar'),
    @('TI_4.png', 0.544, 0.002, 0.511, 1, 'Toni Ismail', 'Benar'),
    @('TI_5.png', 0.892, 0.003, 0.537, 1, 'Toni Ismail', 'Benar'),
    @('RAS_1.png', 0.465, 0.002, 0.498, 1, 'Ridha Ayu Salsabila', 'Benar'),
    @('RAS_2.png', 0.858, 0.003, 0.387, 1, 'Ridha Ayu Salsabila', 'Benar'),
    @('RAS_3.png', 0.363, 0.001, 0.383, 0.857, 'Ridha Ayu Salsabila', 'Benar'),
    @('RAS_4.png', 1.069, 0.004, 0.181, 0.714, 'Ridha Ayu Salsabila', 'Benar'),
    @('RAS_5.png', 0.993, 0.003, 0.502, 1, 'Ridha Ayu Salsabila', 'Benar'),
    @('RR_1.png', 1.34, 0.004, 0.554, 1, 'Rafiqo Rapitasari', 'Benar'),
    @('RR_2.png', 1.278, 0.004, 0.614, 1, 'Rafiqo Rapitasari', 'Benar'),
    @('RR_3.png', 0.976, 0.003, 0.122, 1, 'Rafiqo Rapitasari', 'Benar'),
    @('RR_4.png', 1.225, 0.004, 0.59, 1, 'Rafiqo Rapitasari', 'Benar'),
    @('RR_5.png', 1.34, 0.004, 0.648, 1, 'Rafiqo Rapitasari', 'Benar'),
    @('AR_1.png', 0.678, 0.002, 0.623, 1, 'Arizli Romadhon', 'Benar'),
    @('GA_1.png', 1.526, 0.005, 0.607, 1, 'Gege Ardiyansyah', 'Benar'),
    @('GA_2.png', 0.495, 0.002, 0.308, 0.857, 'Gege Ardiyansyah', 'Benar'),
    @('GA_3.png', 0.577, 0.002, 0.196, 1, 'Gege Ardiyansyah', 'Benar'),
    @('FY_1.png', 0.996, 0.003, 0.403, 1, 'Fanny Yusuf', 'Benar'),
    @('FY_2.png', 1.587, 0.005, 0.522, 1, 'Fanny Yusuf', 'Benar'),
    @('FY_3.png', 1.318, 0.004, 0.539, 0.857, 'Fanny Yusuf', 'Benar'),
    @('FY_4.png', 1.13, 0.004, 0.428, 0.857, 'Fanny Yusuf', 'Benar'),
    @('TO_1.png', 0.631, 0.002, 0.441, 1, 'Tiara Oktavian', 'Benar'),
    @('TO_2.png', 0.902, 0.003, 0.433, 1, 'Tiara Oktavian', 'Benar'),
    @('TO_3.png', 0.67, 0.002, 0.493, 1, 'Tiara Oktavian', 'Benar'),
    @('TO_4.png', 6.848, 0.023, 0.186, 0.857, 'Tiara Oktavian', 'Benar'),
    @('TO_5.png', 6.726, 0.022, 0.137, 0.857, 'Tiara Oktavian', 'Benar'),
    @('TD_1.png', 3.65, 0.012, 0.045, 0.429, 'Tidak Diketahui', 'Benar'),
    @('TD_2.png', 4.798, 0.016, 0.068, 0.429, 'Tidak Diketahui', 'Benar'),
    @('TD_3.png', 0.862, 0.003, 0.211, 0.286, 'Tidak Diketahui', 'Benar'),
    @('TD_4.png', 0.836, 0.003, 0.109, 0.429, 'Tidak Diketahui', 'Benar'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

$ws.Range("A1:G45").NumberFormat = "General"
